# Auto-generated edit script: append rows 655-675 (new NAV data block)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A655").Value = "'2024-08-28"
$ws.Range("C655").Value = 1890.599975585938
$ws.Range("D655").Value = 711.7999877929688
$ws.Range("E655").Value = 78.95999908447266
$ws.Range("F655").Value = 299.9500122070312
$ws.Range("G655").Value = 1323.25
$ws.Range("H655").Value = 29275.72984313965
$ws.Range("I655").Value = 0
$ws.Range("J655").Value = 368.7466846049752

$ws.Range("A656").Value = "'2024-08-29"
$ws.Range("C656").Value = 1893.25
$ws.Range("D656").Value = 700.7999877929688
$ws.Range("E656").Value = 77
$ws.Range("F656").Value = 296.2000122070312
$ws.Range("G656").Value = 1298.650024414062
$ws.Range("H656").Value = 28928.15008544922
$ws.Range("I656").Value = -0.01187262485180639
$ws.Range("J656").Value = 364.3686935533129

$ws.Range("A657").Value = "'2024-08-30"
$ws.Range("C657").Value = 1899.349975585938
$ws.Range("D657").Value = 696.0999755859375
$ws.Range("E657").Value = 75.83999633789062
$ws.Range("F657").Value = 299.2999877929688
$ws.Range("G657").Value = 1302.900024414062
$ws.Range("H657").Value = 28869.96923828125
$ws.Range("I657").Value = -0.002011219071945895
$ws.Range("J657").Value = 363.6358682876185

$ws.Range("A658").Value = "'2024-09-02"
$ws.Range("C658").Value = 1885.400024414062
$ws.Range("D658").Value = 683.5999755859375
$ws.Range("E658").Value = 73.80999755859375
$ws.Range("F658").Value = 296.8999938964844
$ws.Range("G658").Value = 1303.849975585938
$ws.Range("H658").Value = 28487.57955932617
$ws.Range("I658").Value = -0.01324524026329871
$ws.Range("J658").Value = 358.8194238437957

$ws.Range("A659").Value = "'2024-09-03"
$ws.Range("C659").Value = 1901.949951171875
$ws.Range("D659").Value = 689.4000244140625
$ws.Range("E659").Value = 74.47000122070312
$ws.Range("F659").Value = 297.1499938964844
$ws.Range("G659").Value = 1320.25
$ws.Range("H659").Value = 28715.86001586914
$ws.Range("I659").Value = 0.008013332830455756
$ws.Range("J659").Value = 361.6947633130884

$ws.Range("A660").Value = "'2024-09-04"
$ws.Range("C660").Value = 1901.300048828125
$ws.Range("D660").Value = 688.9500122070312
$ws.Range("E660").Value = 74.16000366210938
$ws.Range("F660").Value = 298.9500122070312
$ws.Range("G660").Value = 1327.099975585938
$ws.Range("H660").Value = 28721.03076171875
$ws.Range("I660").Value = 0.0001800658537390795
$ws.Range("J660").Value = 361.7598921894373

$ws.Range("A661").Value = "'2024-09-05"
$ws.Range("C661").Value = 1879.449951171875
$ws.Range("D661").Value = 687.5
$ws.Range("E661").Value = 76
$ws.Range("F661").Value = 290.6000061035156
$ws.Range("G661").Value = 1312.349975585938
$ws.Range("H661").Value = 28602.74978637695
$ws.Range("I661").Value = -0.004118270556621158
$ws.Range("J661").Value = 360.2700670768671

$ws.Range("A662").Value = "'2024-09-06"
$ws.Range("C662").Value = 1872.349975585938
$ws.Range("D662").Value = 673.5499877929688
$ws.Range("E662").Value = 74.72000122070312
$ws.Range("F662").Value = 283.6000061035156
$ws.Range("G662").Value = 1289.699951171875
$ws.Range("H662").Value = 28191.60983276367
$ws.Range("I662").Value = -0.01437414083205038
$ws.Range("J662").Value = 355.091494395132

$ws.Range("A663").Value = "'2024-09-09"
$ws.Range("C663").Value = 1892.400024414062
$ws.Range("D663").Value = 664.1500244140625
$ws.Range("E663").Value = 74.33999633789062
$ws.Range("F663").Value = 281.5499877929688
$ws.Range("G663").Value = 1237.150024414062
$ws.Range("H663").Value = 28036.46997070312
$ws.Range("I663").Value = -0.005503050836077006
$ws.Range("J663").Value = 353.137407850017

$ws.Range("A664").Value = "'2024-09-10"
$ws.Range("C664").Value = 1922.449951171875
$ws.Range("D664").Value = 664.5999755859375
$ws.Range("E664").Value = 78.05000305175781
$ws.Range("F664").Value = 285.75
$ws.Range("G664").Value = 1250.300048828125
$ws.Range("H664").Value = 28561.49984741211
$ws.Range("I664").Value = 0.01872667555001102
$ws.Range("J664").Value = 359.7504975113962

$ws.Range("A665").Value = "'2024-09-11"
$ws.Range("C665").Value = 1957.599975585938
$ws.Range("D665").Value = 689.75
$ws.Range("E665").Value = 81.94999694824219
$ws.Range("F665").Value = 288.0499877929688
$ws.Range("G665").Value = 1237.699951171875
$ws.Range("H665").Value = 29297.64938354492
$ws.Range("I665").Value = 0.0257741904334731
$ws.Range("J665").Value = 369.0227753427916

$ws.Range("A666").Value = "'2024-09-12"
$ws.Range("C666").Value = 1996.400024414062
$ws.Range("D666").Value = 729.1500244140625
$ws.Range("E666").Value = 81.69999694824219
$ws.Range("F666").Value = 291.7000122070312
$ws.Range("G666").Value = 1237.300048828125
$ws.Range("H666").Value = 29912.80038452148
$ws.Range("I666").Value = 0.02099659917843318
$ws.Range("J666").Value = 376.7709986443772

$ws.Range("A667").Value = "'2024-09-13"
$ws.Range("C667").Value = 1988.050048828125
$ws.Range("D667").Value = 713.7000122070312
$ws.Range("E667").Value = 83.11000061035156
$ws.Range("F667").Value = 289.9500122070312
$ws.Range("G667").Value = 1241.5
$ws.Range("H667").Value = 29812.18057250977
$ws.Range("I667").Value = -0.003363771051799782
$ws.Range("J667").Value = 375.5036272659796

$ws.Range("A668").Value = "'2024-09-16"
$ws.Range("C668").Value = 1989.900024414062
$ws.Range("D668").Value = 714.2000122070312
$ws.Range("E668").Value = 84.69999694824219
$ws.Range("F668").Value = 290.3999938964844
$ws.Range("G668").Value = 1226.599975585938
$ws.Range("H668").Value = 29926.49987792969
$ws.Range("I668").Value = 0.003834650911961043
$ws.Range("J668").Value = 376.9435525927198

$ws.Range("A669").Value = "'2024-09-17"
$ws.Range("C669").Value = 2006.550048828125
$ws.Range("D669").Value = 731.0999755859375
$ws.Range("E669").Value = 82
$ws.Range("F669").Value = 284.2999877929688
$ws.Range("G669").Value = 1193.800048828125
$ws.Range("H669").Value = 29823.24993896484
$ws.Range("I669").Value = -0.003450117433913109
$ws.Range("J669").Value = 375.6430530703184

$ws.Range("A670").Value = "'2024-09-18"
$ws.Range("C670").Value = 1987.800048828125
$ws.Range("D670").Value = 743.25
$ws.Range("E670").Value = 80.81999969482422
$ws.Range("F670").Value = 282.8500061035156
$ws.Range("G670").Value = 1166.400024414062
$ws.Range("H670").Value = 29685.31034851074
$ws.Range("I670").Value = -0.004625236710834788
$ws.Range("J670").Value = 373.9056150310876

$ws.Range("A671").Value = "'2024-09-19"
$ws.Range("C671").Value = 1998.599975585938
$ws.Range("D671").Value = 735.9500122070312
$ws.Range("E671").Value = 80.97000122070312
$ws.Range("F671").Value = 272.7000122070312
$ws.Range("G671").Value = 1121.300048828125
$ws.Range("H671").Value = 29455.8603515625
$ws.Range("I671").Value = -0.0077294120982553
$ws.Range("J671").Value = 371.0155444466607

$ws.Range("A672").Value = "'2024-09-20"
$ws.Range("C672").Value = 2048.10009765625
$ws.Range("D672").Value = 746.5
$ws.Range("E672").Value = 83.44999694824219
$ws.Range("F672").Value = 277.3500061035156
$ws.Range("G672").Value = 1149.400024414062
$ws.Range("H672").Value = 30118.95037841797
$ws.Range("I672").Value = 0.02251131078642199
$ws.Range("J672").Value = 379.3675906742931

$ws.Range("A673").Value = "'2024-09-23"
$ws.Range("C673").Value = 2082.39990234375
$ws.Range("D673").Value = 773.9500122070312
$ws.Range("E673").Value = 82.88999938964844
$ws.Range("F673").Value = 286.2999877929688
$ws.Range("G673").Value = 1162.75
$ws.Range("H673").Value = 30664.31942749023
$ws.Range("I673").Value = 0.01810717313253569
$ws.Range("J673").Value = 386.2368653195055

$ws.Range("A674").Value = "'2024-09-24"
$ws.Range("C674").Value = 2068.14990234375
$ws.Range("D674").Value = 781.8499755859375
$ws.Range("E674").Value = 83.79000091552734
$ws.Range("F674").Value = 291.7999877929688
$ws.Range("G674").Value = 1141.199951171875
$ws.Range("H674").Value = 30770.6690826416
$ws.Range("I674").Value = 0.003468188994144963
$ws.Range("J674").Value = 387.5764077649396

$ws.Range("A675").Value = "'2024-09-25"
$ws.Range("C675").Value = 2061.60009765625
$ws.Range("D675").Value = 775.8499755859375
$ws.Range("E675").Value = 82.95999908447266
$ws.Range("F675").Value = 289.8500061035156
$ws.Range("G675").Value = 1118.449951171875
$ws.Range("H675").Value = 30542.33015441895
$ws.Range("I675").Value = -0.007420668286718119
$ws.Range("J675").Value = 384.7003318071582

